# Applies the "Add all my files" commit:
#  - Inserts two new worksheets, "ProfileData" and "Search", between
#    "SignIn" and "ShareSkill".
#  - Adds a third row of data (an extra saved login) to the "SignIn" sheet.
#  - Updates two date values on the "ShareSkill" sheet (keeping their
#    original date-format style) and moves its active-cell selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. SignIn sheet: append a 3rd saved-login row (B3/C3 only, no A3).
# ---------------------------------------------------------------------
$signIn = $wb.Worksheets.Item("SignIn")
$signIn.Range("B3").Value = "mvpstudio.qa@gmail.com"
$signIn.Range("C3").Value = "SydneyQa2019"
$signIn.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Insert "ProfileData" right after "SignIn".
# ---------------------------------------------------------------------
$profileData = $wb.Worksheets.Add($null, $signIn)
$profileData.Name = "ProfileData"

# Row 1 header, Row 2 data -- written in the same left-to-right,
# top-to-bottom-ish order the strings were authored in.
$profileData.Range("B1").Value = "Lname"
$profileData.Range("A2").Value = "MVP"
$profileData.Range("B2").Value = "Studio"
$profileData.Range("A1").Value = "Fname"
$profileData.Range("C1").Value = "Availability"
$profileData.Range("D1").Value = "Hours"
$profileData.Range("E1").Value = "EarnTarget"
$profileData.Range("C2").Value = "Full Time"
$profileData.Range("D2").Value = "More than 30hours a week"
$profileData.Range("E2").Value = "More than $1000 per month"
$profileData.Range("F1").Value = "DescriptionData"
$profileData.Range("F2").Value = "I am a QA Engineer"

$profileData.Columns.Item(3).ColumnWidth = 13.25
$profileData.Columns.Item(4).ColumnWidth = 27.25
$profileData.Columns.Item(5).ColumnWidth = 29.42
$profileData.Columns.Item(6).ColumnWidth = 25.59

$profileData.Range("F9").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Insert "Search" right after "ProfileData".
# ---------------------------------------------------------------------
$search = $wb.Worksheets.Add($null, $profileData)
$search.Name = "Search"

$search.Range("B1").Value = "User"
$search.Range("B2").Value = "Priyanka Singh"
$search.Range("A2").Value = "I design beautiful logos"
$search.Range("A1").Value = "Title"

$search.Columns.Item(1).ColumnWidth = 37.92
$search.Columns.Item(2).ColumnWidth = 14.92

$search.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. ShareSkill sheet: bump the two Startdate values to 44016
#    (04/07/2020) while preserving their original date-format style,
#    and move the active-cell selection to H6.
# ---------------------------------------------------------------------
$shareSkill = $wb.Worksheets.Item("ShareSkill")

$shareSkill.Range("H2").Copy() | Out-Null
$shareSkill.Range("Z1").PasteSpecial(-4122) | Out-Null
$shareSkill.Range("H2").Value = 44016
$shareSkill.Range("Z1").Copy() | Out-Null
$shareSkill.Range("H2").PasteSpecial(-4122) | Out-Null
$shareSkill.Range("Z1").Clear() | Out-Null

$shareSkill.Range("H3").Copy() | Out-Null
$shareSkill.Range("Z1").PasteSpecial(-4122) | Out-Null
$shareSkill.Range("H3").Value = 44016
$shareSkill.Range("Z1").Copy() | Out-Null
$shareSkill.Range("H3").PasteSpecial(-4122) | Out-Null
$shareSkill.Range("Z1").Clear() | Out-Null

$shareSkill.Select() | Out-Null
$shareSkill.Range("H6").Select() | Out-Null
